# Update statistics table with values from the optimized DOE run (5000 cases)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Evaporator Temperature
$ws.Range("C2").Value = 14.95002
$ws.Range("D2").Value = 8.662638643118635
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 7.475
$ws.Range("G2").Value = 14.95
$ws.Range("H2").Value = 22.425
$ws.Range("I2").Value = 30

# Row 3: Condenser Temperature
$ws.Range("C3").Value = 49.95002000000001
$ws.Range("D3").Value = 8.662638643118635
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 42.475
$ws.Range("G3").Value = 49.95
$ws.Range("H3").Value = 57.42500000000001
$ws.Range("I3").Value = 65

# Row 4: Adiabatic Efficiency
$ws.Range("C4").Value = 74.95001999999999
$ws.Range("D4").Value = 2.886930354714951
$ws.Range("E4").Value = 70
$ws.Range("F4").Value = 72.47499999999999
$ws.Range("G4").Value = 74.95
$ws.Range("H4").Value = 77.42500000000001
$ws.Range("I4").Value = 80

# Row 5: Compressor Energy
$ws.Range("C5").Value = 883.9745256500373
$ws.Range("D5").Value = 548.9691248772526
$ws.Range("E5").Value = 51.06911025303199
$ws.Range("F5").Value = 473.3640293741739
$ws.Range("G5").Value = 757.487449361784
$ws.Range("H5").Value = 1168.897416824464
$ws.Range("I5").Value = 3426.78767522376

# Row 6: Electric Current
$ws.Range("C6").Value = 4.018066025681987
$ws.Range("D6").Value = 2.495314203987512
$ws.Range("E6").Value = 0.2321323193319636
$ws.Range("F6").Value = 2.151654678973517
$ws.Range("G6").Value = 3.443124769826291
$ws.Range("H6").Value = 5.313170076474837
$ws.Range("I6").Value = 15.57630761465346

# Row 7: Discharge Temperature
$ws.Range("C7").Value = 69.97073485711803
$ws.Range("D7").Value = 14.29305141868105
$ws.Range("E7").Value = 38.06536296353289
$ws.Range("F7").Value = 58.76314971746178
$ws.Range("G7").Value = 69.9688202484181
$ws.Range("H7").Value = 81.61145300906441
$ws.Range("I7").Value = 103.6520532969717

# Row 8: Refrigerant Mass Flow
$ws.Range("C8").Value = 1.620936784243737
$ws.Range("D8").Value = 0.6183514256259423
$ws.Range("E8").Value = 0.5312047499164733
$ws.Range("F8").Value = 1.109644836230542
$ws.Range("G8").Value = 1.588888018984051
$ws.Range("H8").Value = 2.066226127260754
$ws.Range("I8").Value = 3.386818343590615

# Row 9: Capacity
$ws.Range("C9").Value = 12450.02
$ws.Range("D9").Value = 4331.276023298215
$ws.Range("E9").Value = 5000
$ws.Range("F9").Value = 8700
$ws.Range("G9").Value = 12450
$ws.Range("H9").Value = 16200
$ws.Range("I9").Value = 20000

$wb.Save()
